$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("variable_mapping")

# Insert a new row above the existing "Naver" row (currently row 5) to host
# the new "Coupang" platform mapping. This shifts Naver down to row 6.
$ws.Rows.Item(5).Insert()
$ws.Rows.Item(5).RowHeight = 17.25

# Nudge J5 so it resolves to a clean "blank" cell state (matches the blank
# K5/J6-style neighbour cells instead of inheriting a stale quote-prefix style).
$ws.Range("J5").Value = "x"
$ws.Range("J5").ClearContents()

# Fill in the Coupang header/mapping row.
$ws.Range("A5").Value = "Coupang"
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "주문번호"
$ws.Range("D5").Value = "등록상품명"
$ws.Range("E5").Value = "등록옵션명"
$ws.Range("F5").Value = "구매수(수량)"
$ws.Range("G5").Value = "수취인이름"
$ws.Range("H5").Value = "우편번호"
$ws.Range("I5").Value = "수취인 주소"
$ws.Range("L5").Value = "수취인전화번호"
$ws.Range("M5").Value = "구매자전화번호"
$ws.Range("N5").Value = "배송메세지"
$ws.Range("O5").Value = "노출상품ID"
